# Apply crypto price/volume updates per the commit diff.
# Data rows for this sheet run 2..51 (Coin=B, Link=C, Price=D, Volume(1h)=E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price (D) values are purely numeric-looking text (e.g. "226.99") that
# Excel would otherwise auto-convert to a number on assignment. Force those
# specific cells to remain text, matching the original inlineStr cells, then
# restore the default "Normal" style so no stray formatting is introduced.
$textForceCells = @("D5", "D6", "D8", "D9", "D10", "D13", "D14", "D15", "D16", "D20", "D22", "D25", "D26", "D27", "D28", "D29", "D32", "D38", "D39", "D41", "D43", "D44", "D45", "D48", "D49", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '37.162.10'
$ws.Range("E2").Value = '  -1.61%  '
$ws.Range("D3").Value = '2.026.27'
$ws.Range("E3").Value = '  -2.96%  '
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '226.99'
$ws.Range("E5").Value = '  -3.02%  '
$ws.Range("D6").Value = '0.608'
$ws.Range("E6").Value = '  -2.81%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").Value = '55.40'
$ws.Range("E8").Value = '  -4.92%  '
$ws.Range("D9").Value = '0.382'
$ws.Range("E9").Value = '  -3.00%  '
$ws.Range("D10").Value = '0.0792'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("E11").Value = '  -2.64%  '
$ws.Range("D12").Value = '2.327.69'
$ws.Range("E12").Value = '  -2.84%  '
$ws.Range("D13").Value = '14.32'
$ws.Range("E13").Value = '  -5.85%  '
$ws.Range("D14").Value = '20.35'
$ws.Range("E14").Value = '  -3.96%  '
$ws.Range("D15").Value = '0.744'
$ws.Range("E15").Value = '  -4.36%  '
$ws.Range("D16").Value = '5.18'
$ws.Range("E16").Value = '  -3.26%  '
$ws.Range("D17").Value = '2.035.72'
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = '37.013.80'
$ws.Range("E18").Value = '  -1.99%  '
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").Value = '68.95'
$ws.Range("E20").Value = '  -2.94%  '
$ws.Range("D21").Value = '0.0₃0835'
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("D22").Value = '223.58'
$ws.Range("E22").Value = '  -2.70%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '2.26'
$ws.Range("E25").Value = '  -5.83%  '
$ws.Range("D26").Value = '9.38'
$ws.Range("E26").Value = '  -3.88%  '
$ws.Range("D27").Value = '167.88'
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("D28").Value = '0.125'
$ws.Range("E28").Value = '  -6.74%  '
$ws.Range("D29").Value = '18.74'
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("E30").Value = '  -3.65%  '
$ws.Range("E31").Value = '  -4.04%  '
$ws.Range("D32").Value = '4.47'
$ws.Range("E32").Value = '  -4.83%  '
$ws.Range("E33").Value = '  -4.35%  '
$ws.Range("E34").Value = '  -3.32%  '
$ws.Range("E35").Value = '  -5.45%  '
$ws.Range("E36").Value = '  +0.38%  '
$ws.Range("E37").Value = '  +0.02%  '
$ws.Range("D38").Value = '3.15'
$ws.Range("E38").Value = '  -5.13%  '
$ws.Range("D39").Value = '5.32'
$ws.Range("E39").Value = '  -1.07%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.489.47'
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.0218'
$ws.Range("E41").Value = '  -7.85%  '
$ws.Range("E42").Value = '  -2.12%  '
$ws.Range("D43").Value = '0.0930'
$ws.Range("E43").Value = '  -4.14%  '
$ws.Range("D44").Value = '95.03'
$ws.Range("E44").Value = '  -6.02%  '
$ws.Range("D45").Value = '16.47'
$ws.Range("E45").Value = '  -1.26%  '
$ws.Range("E46").Value = '  -6.14%  '
$ws.Range("E47").Value = '  -5.27%  '
$ws.Range("D48").Value = '7.12'
$ws.Range("E48").Value = '  -1.20%  '
$ws.Range("D49").Value = '2.90'
$ws.Range("E49").Value = '  -1.86%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.218.85'
$ws.Range("E50").Value = '  -2.65%  '
$ws.Range("B51").Value = 'FTXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D51").Value = '3.66'
$ws.Range("E51").Value = '  -11.03%  '

foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}

